$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

$ws1.Range("H86").Value = 57711092
$ws1.Range("I86").Value = 72223840
$ws1.Range("J86").Value = 9335259
$ws1.Range("K86").Value = 72223840
$ws1.Range("L86").Value = 9335259
$ws1.Range("M86").Value = -72222717
$ws1.Range("N86").Value = -9337505

$ws1.Range("H89").Value = 57711092
$ws1.Range("I89").Value = 72223840
$ws1.Range("J89").Value = 9335259
$ws1.Range("K89").Value = 361119200
$ws1.Range("L89").Value = 46676295
$ws1.Range("M89").Value = -361113584
$ws1.Range("N89").Value = -46687527

$ws1.Range("H98").Value = 32260884
$ws1.Range("I98").Value = 34485524
$ws1.Range("J98").Value = 3598
$ws1.Range("K98").Value = 34485524
$ws1.Range("L98").Value = 3598
$ws1.Range("M98").Value = -34484026
$ws1.Range("N98").Value = -6594

$ws1.Range("H122").Value = 32260884
$ws1.Range("I122").Value = 34485524
$ws1.Range("J122").Value = 3598
$ws1.Range("K122").Value = 103456572
$ws1.Range("L122").Value = 10794
$ws1.Range("M122").Value = -103454122
$ws1.Range("N122").Value = -15694

$ws1.Range("H125").Value = 29413172
$ws1.Range("J125").Value = 3997
$ws1.Range("L125").Value = 35973
$ws1.Range("N125").Value = -40893

$ws1.Range("H129").Value = 1793.5
$ws1.Range("I129").Value = 1793.5
$ws1.Range("K129").Value = 5380.5
$ws1.Range("M129").Value = -380.5

$ws1.Range("H132").Value = 2686.7036
$ws1.Range("I132").Value = 2590.5686
$ws1.Range("J132").Value = 4321
$ws1.Range("K132").Value = 7771.7058
$ws1.Range("L132").Value = 12963
$ws1.Range("M132").Value = -5241.7058
$ws1.Range("N132").Value = -18023

$ws1.Range("H135").Value = 715094.8
$ws1.Range("I135").Value = 909520.6
$ws1.Range("J135").Value = 2200
$ws1.Range("K135").Value = 8185685.399999999
$ws1.Range("L135").Value = 19800
$ws1.Range("M135").Value = -8183150.399999999
$ws1.Range("N135").Value = -24870

$ws1.Range("H137").Value = 3738.889
$ws1.Range("I137").Value = 5600.25
$ws1.Range("K137").Value = 16800.75
$ws1.Range("M137").Value = -14250.75

$ws1.Range("H138").Value = 5500.5317
$ws1.Range("I138").Value = 1230
$ws1.Range("K138").Value = 3690
$ws1.Range("M138").Value = 1450

$ws1.Range("H141").Value = 3999
$ws1.Range("I141").Value = 3747.5
$ws1.Range("J141").Value = 4250.5
$ws1.Range("K141").Value = 11242.5
$ws1.Range("L141").Value = 12751.5
$ws1.Range("M141").Value = -6062.5
$ws1.Range("N141").Value = -23111.5

$ws2.Range("H63").Value = 1765
$ws2.Range("I63").Value = 1765
$ws2.Range("J63").Value = 0
$ws2.Range("K63").Value = 1765
$ws2.Range("L63").Value = 0
$ws2.Range("M63").Value = -1079
$ws2.Range("N63").ClearContents()

$ws2.Range("H66").Value = 1765
$ws2.Range("I66").Value = 1765
$ws2.Range("J66").Value = 0
$ws2.Range("K66").Value = 8825
$ws2.Range("L66").Value = 0
$ws2.Range("M66").Value = -5393
$ws2.Range("N66").ClearContents()

$ws2.Range("H74").Value = 82876.64999999999
$ws2.Range("I74").Value = 147026.55
$ws2.Range("J74").Value = 4471.222
$ws2.Range("K74").Value = 147026.55
$ws2.Range("L74").Value = 4471.222
$ws2.Range("M74").Value = -146152.55
$ws2.Range("N74").Value = -6219.222

$ws2.Range("H77").Value = 82876.64999999999
$ws2.Range("I77").Value = 147026.55
$ws2.Range("J77").Value = 4471.222
$ws2.Range("K77").Value = 735132.75
$ws2.Range("L77").Value = 22356.11
$ws2.Range("M77").Value = -730764.75
$ws2.Range("N77").Value = -31092.11

$ws2.Range("H88").Value = 1918.4117
$ws2.Range("I88").Value = 1695.125
$ws2.Range("J88").Value = 2116.889
$ws2.Range("K88").Value = 1695.125
$ws2.Range("L88").Value = 2116.889
$ws2.Range("M88").Value = -1289.125
$ws2.Range("N88").Value = -2928.889

$ws2.Range("H91").Value = 1918.4117
$ws2.Range("I91").Value = 1695.125
$ws2.Range("J91").Value = 2116.889
$ws2.Range("K91").Value = 1695.125
$ws2.Range("L91").Value = 2116.889
$ws2.Range("M91").Value = -291.125
$ws2.Range("N91").Value = -4924.889

$ws2.Range("H102").Value = 4744.727
$ws2.Range("I102").Value = 4600.4287
$ws2.Range("K102").Value = 4600.4287
$ws2.Range("M102").Value = -2978.4287

$ws2.Range("H122").Value = 14768.611
$ws2.Range("I122").Value = 18064.691
$ws2.Range("J122").Value = 6198.8
$ws2.Range("K122").Value = 54194.073
$ws2.Range("L122").Value = 18596.4
$ws2.Range("M122").Value = -51744.073
$ws2.Range("N122").Value = -23496.4

$ws3.Range("H82").Value = 5332.6665
$ws3.Range("I82").Value = 5332.6665
$ws3.Range("K82").Value = 5332.6665
$ws3.Range("M82").Value = -4949.6665

$ws3.Range("H85").Value = 5332.6665
$ws3.Range("I85").Value = 5332.6665
$ws3.Range("K85").Value = 5332.6665
$ws3.Range("M85").Value = -4006.6665

$ws4.Range("H31").Value = 5912.9717
$ws4.Range("I31").Value = 2512.4546
$ws4.Range("J31").Value = 11454.556
$ws4.Range("K31").Value = 2512.4546
$ws4.Range("L31").Value = 11454.556
$ws4.Range("M31").Value = -2217.4546
$ws4.Range("N31").Value = -12044.556

$ws4.Range("H34").Value = 5912.9717
$ws4.Range("I34").Value = 2512.4546
$ws4.Range("J34").Value = 11454.556
$ws4.Range("K34").Value = 2512.4546
$ws4.Range("L34").Value = 11454.556
$ws4.Range("M34").Value = -2310.4546
$ws4.Range("N34").Value = -11858.556

$ws4.Range("H62").Value = 6810.5713
$ws4.Range("I62").Value = 5534.8
$ws4.Range("K62").Value = 5534.8
$ws4.Range("M62").Value = -4910.8

$ws4.Range("H65").Value = 6810.5713
$ws4.Range("I65").Value = 5534.8
$ws4.Range("K65").Value = 27674
$ws4.Range("M65").Value = -24554

$ws4.Range("H68").Value = 79992.664
$ws4.Range("J68").Value = 79992.664
$ws4.Range("L68").Value = 79992.664
$ws4.Range("N68").Value = -81490.664

$ws4.Range("H71").Value = 79992.664
$ws4.Range("J71").Value = 79992.664
$ws4.Range("L71").Value = 239977.992
$ws4.Range("N71").Value = -247465.992

$ws4.Range("H74").Value = 333400000
$ws4.Range("J74").Value = 99998
$ws4.Range("L74").Value = 99998
$ws4.Range("N74").Value = -101746

$ws4.Range("H77").Value = 333400000
$ws4.Range("J77").Value = 99998
$ws4.Range("L77").Value = 299994
$ws4.Range("N77").Value = -308730

$ws4.Range("H107").Value = 2234.5881
$ws4.Range("J107").Value = 2295.818
$ws4.Range("L107").Value = 2295.818
$ws4.Range("N107").Value = -6135.818

$ws5.Range("H132").Value = 10430.883
$ws5.Range("J132").Value = 15814.125
$ws5.Range("L132").Value = 142327.125
$ws5.Range("N132").Value = -147387.125

$ws5.Range("H141").Value = 15750

$ws6.Range("H57").Value = 60001.383
$ws6.Range("J57").Value = 64996.918
$ws6.Range("L57").Value = 64996.918
$ws6.Range("N57").Value = -66636.91800000001

$ws6.Range("H70").Value = 7482.385
$ws6.Range("I70").Value = 6070.1577
$ws6.Range("J70").Value = 11315.571
$ws6.Range("K70").Value = 6070.1577
$ws6.Range("L70").Value = 11315.571
$ws6.Range("M70").Value = -5800.1577
$ws6.Range("N70").Value = -11855.571

$ws6.Range("H73").Value = 7482.385
$ws6.Range("I73").Value = 6070.1577
$ws6.Range("J73").Value = 11315.571
$ws6.Range("K73").Value = 6070.1577
$ws6.Range("L73").Value = 11315.571
$ws6.Range("M73").Value = -5134.1577
$ws6.Range("N73").Value = -13187.571

$ws7.Range("H61").Value = 5566.8887
$ws7.Range("I61").Value = 3843.4285
$ws7.Range("K61").Value = 3843.4285
$ws7.Range("M61").Value = -3641.4285

$ws7.Range("H113").Value = 5566.8887
$ws7.Range("I113").Value = 3843.4285
$ws7.Range("K113").Value = 3843.4285
$ws7.Range("M113").Value = -1673.4285

$ws7.Range("H122").Value = 3866.5
$ws7.Range("I122").Value = 3414.9
$ws7.Range("K122").Value = 10244.7
$ws7.Range("M122").Value = -7794.700000000001

$ws7.Range("H132").Value = 12506538
$ws7.Range("I132").Value = 19234420
$ws7.Range("K132").Value = 57703260
$ws7.Range("M132").Value = -57700730

$ws8.Range("H96").Value = 1810.7142
$ws8.Range("I96").Value = 1696
$ws8.Range("K96").Value = 1696
$ws8.Range("M96").Value = -323

$ws8.Range("H122").Value = 6594.1113
$ws8.Range("I122").Value = 5083.5
$ws8.Range("K122").Value = 15250.5
$ws8.Range("M122").Value = -12800.5

$ws8.Range("H132").Value = 11644112
$ws8.Range("I132").Value = 15155942
$ws8.Range("K132").Value = 45467826
$ws8.Range("M132").Value = -45465296

$ws8.Range("H141").Value = 0
$ws8.Range("J141").Value = 0
$ws8.Range("L141").Value = 0
$ws8.Range("N141").ClearContents()
